$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4, column F (acc values)
$ws.Range("F2").Value = 0.8744787322768974
$ws.Range("F3").Value = 0.8836530442035029
$ws.Range("F4").Value = 0.8949124270225187

# New rows data: model_type, learning_rate, hidden_size, batch_size, pooling_style, acc
$data = @(
    @("gated_cnn", 0.001, 128, 256, "max", 0.8878231859883235),
    @("gated_cnn", 0.001, 256, 64, "avg", 0.8669724770642202),
    @("gated_cnn", 0.001, 256, 64, "max", 0.8811509591326105),
    @("gated_cnn", 0.001, 256, 256, "avg", 0.8811509591326105),
    @("gated_cnn", 0.001, 256, 256, "max", 0.8886572143452878),
    @("gated_cnn", 0.0001, 128, 64, "avg", 0.8786488740617181),
    @("gated_cnn", 0.0001, 128, 64, "max", 0.8928273561301084),
    @("gated_cnn", 0.0001, 128, 256, "avg", 0.8432026688907422),
    @("gated_cnn", 0.0001, 128, 256, "max", 0.872393661384487),
    @("fast_text", 0.001, 128, 64, "avg", 0.8915763135946623),
    @("lstm", 0.001, 128, 64, "avg", 0.8798999165971643)
)

$row = 5
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $row++
}
